# Updated cryptos list — applies the new Price (column D) and Volume(1h)
# (column E) figures for each coin row, matching the latest scrape.
#
# Column D values that look like plain numbers (e.g. "0.7025") are forced
# back to text via NumberFormat "@" before assignment, matching the
# original inline-string ("Text"-typed) cells -- otherwise Excel's COM
# layer auto-coerces a bare numeric-looking string into a real number.
# Values that already contain non-numeric punctuation (e.g. "29.327.31",
# which has two dots) or the padded "  +0.27%  " percentage strings stay
# text on their own because Excel can't parse them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "29.327.31" },
    @{ Cell = "E2"; Value = "  +0.27%  " },
    @{ Cell = "D3"; Value = "1.861.35" },
    @{ Cell = "E3"; Value = "  +0.18%  " },
    @{ Cell = "E4"; Value = "  +0.07%  " },
    @{ Cell = "D5"; Value = "0.7025" },
    @{ Cell = "E5"; Value = "  +0.38%  " },
    @{ Cell = "E6"; Value = "  +0.13%  " },
    @{ Cell = "D8"; Value = "0.07839" },
    @{ Cell = "E8"; Value = "  -1.86%  " },
    @{ Cell = "D9"; Value = "0.3049" },
    @{ Cell = "E9"; Value = "  -0.05%  " },
    @{ Cell = "D10"; Value = "24.72" },
    @{ Cell = "E10"; Value = "  +6.15%  " },
    @{ Cell = "D11"; Value = "0.08150" },
    @{ Cell = "E11"; Value = "  -0.48%  " },
    @{ Cell = "D12"; Value = "1.868.05" },
    @{ Cell = "E12"; Value = "  +0.31%  " },
    @{ Cell = "E13"; Value = "  +0.85%  " },
    @{ Cell = "D14"; Value = "0.7135" },
    @{ Cell = "E14"; Value = "  -0.60%  " },
    @{ Cell = "D15"; Value = "89.15" },
    @{ Cell = "E15"; Value = "  +0.13%  " },
    @{ Cell = "D16"; Value = "29.361.66" },
    @{ Cell = "E16"; Value = "  +0.38%  " },
    @{ Cell = "D17"; Value = "5.804" },
    @{ Cell = "E17"; Value = "  +0.82%  " },
    @{ Cell = "D18"; Value = "0.000007779" },
    @{ Cell = "E18"; Value = "  -0.13%  " },
    @{ Cell = "D19"; Value = "239.17" },
    @{ Cell = "E19"; Value = "  +1.04%  " },
    @{ Cell = "D20"; Value = "13.17" },
    @{ Cell = "E20"; Value = "  -1.22%  " },
    @{ Cell = "D21"; Value = "2.125.31" },
    @{ Cell = "E21"; Value = "  +1.14%  " },
    @{ Cell = "E22"; Value = "  +0.11%  " },
    @{ Cell = "E23"; Value = "  +0.13%  " },
    @{ Cell = "D24"; Value = "7.510" },
    @{ Cell = "E24"; Value = "  +0.90%  " },
    @{ Cell = "D25"; Value = "162.66" },
    @{ Cell = "E25"; Value = "  +0.54%  " },
    @{ Cell = "D26"; Value = "8.897" },
    @{ Cell = "E26"; Value = "  -1.06%  " },
    @{ Cell = "D27"; Value = "0.1423" },
    @{ Cell = "E27"; Value = "  -2.28%  " },
    @{ Cell = "D28"; Value = "18.06" },
    @{ Cell = "E28"; Value = "  -0.03%  " },
    @{ Cell = "D29"; Value = "1.902" },
    @{ Cell = "E29"; Value = "  -5.05%  " },
    @{ Cell = "D30"; Value = "1.376" },
    @{ Cell = "E30"; Value = "  -4.23%  " },
    @{ Cell = "D31"; Value = "1.471" },
    @{ Cell = "E31"; Value = "  -0.82%  " },
    @{ Cell = "D32"; Value = "4.295" },
    @{ Cell = "E32"; Value = "  -2.59%  " },
    @{ Cell = "D33"; Value = "4.036" },
    @{ Cell = "E33"; Value = "  -0.34%  " },
    @{ Cell = "D34"; Value = "0.05170" },
    @{ Cell = "E34"; Value = "  -0.94%  " },
    @{ Cell = "D35"; Value = "1.178" },
    @{ Cell = "E35"; Value = "  +0.72%  " },
    @{ Cell = "D36"; Value = "0.7054" },
    @{ Cell = "E36"; Value = "  -0.06%  " },
    @{ Cell = "D37"; Value = "0.9981" },
    @{ Cell = "E37"; Value = "  -0.25%  " },
    @{ Cell = "D38"; Value = "2.677" },
    @{ Cell = "E38"; Value = "  +0.51%  " },
    @{ Cell = "E39"; Value = "  +0.00%  " },
    @{ Cell = "D40"; Value = "2.694" },
    @{ Cell = "D41"; Value = "1.172.87" },
    @{ Cell = "E41"; Value = "  +2.67%  " },
    @{ Cell = "D42"; Value = "0.9183" },
    @{ Cell = "E42"; Value = "  -0.44%  " },
    @{ Cell = "D43"; Value = "6.018" },
    @{ Cell = "E43"; Value = "  +1.59%  " },
    @{ Cell = "D44"; Value = "71.55" },
    @{ Cell = "E44"; Value = "  +1.16%  " },
    @{ Cell = "D45"; Value = "0.4247" },
    @{ Cell = "E45"; Value = "  -0.64%  " },
    @{ Cell = "E46"; Value = "  +0.01%  " },
    @{ Cell = "D47"; Value = "101.81" },
    @{ Cell = "E47"; Value = "  -1.46%  " },
    @{ Cell = "D48"; Value = "0.5351" },
    @{ Cell = "E48"; Value = "  -1.40%  " },
    @{ Cell = "E49"; Value = "  -2.50%  " },
    @{ Cell = "D50"; Value = "9.145" },
    @{ Cell = "E50"; Value = "  -0.30%  " },
    @{ Cell = "D51"; Value = "6.953" },
    @{ Cell = "E51"; Value = "  -0.29%  " }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $v = $u.Value
    $isNumericLooking = $v -match '^[+-]?\d+(\.\d+)?$'
    if ($isNumericLooking) {
        # Force text storage so cells like Price stay strings (as in the
        # source data) instead of being auto-converted to a Double.
        $range.NumberFormat = "@"
    }
    $range.Value = $v
}

